# ui_parameter.xlsx: "image included and ui changed"
#
# The sheet is a sequence of 6-row blocks, each describing one UI image
# asset (Name / Filename / Position_X / Position_Y / Size_X / Size_Y).
# This edit swaps which images/identifiers are used in several blocks,
# tweaks their numeric placement values, adds a brand-new "PLAYER_TIMER"
# block (rows 31-36) replacing the old PLAYER_ATTACK_TYPE block there,
# and moves PLAYER_ATTACK_TYPE/icon_basic.png up into the block that used
# to be ITEM_BAR (rows 25-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String-valued cells first, in the exact order the new names were
# introduced upstream, so the shared-strings table comes out in the same
# order as the authored workbook. ----------------------------------------
$ws.Range("B25").Value = "PLAYER_ATTACK_TYPE"
$ws.Range("B26").Value = "icon_basic.png"
$ws.Range("B31").Value = "PLAYER_TIMER"
$ws.Range("B2").Value = "boss_hp_back.png"
$ws.Range("B8").Value = "boss_hp_front.png"
$ws.Range("B32").Value = "timer2.png"
$ws.Range("B14").Value = "player_hp_back.png"
$ws.Range("B20").Value = "player_hp_front.png"

# --- Block: rows 1-6 (BOSS_HPBAR_BACK) ---------------------------------
$ws.Range("B3").Value = 1423
$ws.Range("B4").Value = 32
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0

# --- Block: rows 7-12 (BOSS_HPBAR_FRONT) -------------------------------
$ws.Range("B9").Value = 1442
$ws.Range("B10").Value = 37
$ws.Range("B11").Value = 426
$ws.Range("B12").Value = 23

# --- Block: rows 13-18 (PLAYER_HPBAR_BACK) -----------------------------
$ws.Range("B15").Value = 125
$ws.Range("B16").Value = 1008
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0

# --- Block: rows 19-24 (PLAYER_HPBAR_FRONT) ----------------------------
$ws.Range("B21").Value = 139
$ws.Range("B22").Value = 1012.5
$ws.Range("B23").Value = 282
$ws.Range("B24").Value = 25

# --- Block: rows 25-30 -------------------------------------------------
# Used to be ITEM_BAR/item_bar.png; now holds PLAYER_ATTACK_TYPE/icon_basic.png
# (names already set above). The trailing "use original image size" comment
# moves from C28/C29 down to C29/C30.
$ws.Range("B27").Value = 30
$ws.Range("B28").Value = 771
$ws.Range("C28").ClearContents()
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = "이미지 원래 사이즈를 사용하고 싶으면 0"

# --- Block: rows 31-36 --------------------------------------------------
# Used to be PLAYER_ATTACK_TYPE/icon_basic.png; now a new PLAYER_TIMER block
# (names already set above). The trailing comment column (C) is dropped
# entirely for this block (cells removed, not just blanked), and the block
# shrinks to columns A:B.
$ws.Range("B33").Value = 8
$ws.Range("B34").Value = 877
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("C31:C36").Clear()

# --- View state: scroll position + active selection --------------------
$ws.Range("D10").Select()
